$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record arrived for "Vega Modelo de Temuco - Alcachofa".
# It belongs between the existing rows 127 and 128 (by date order used in
# the sheet), so insert a fresh row at 128 - this shifts the old rows
# 128-147 down to 129-148 - and fill it with the new record's data.
$ws.Rows(128).Insert()

$ws.Range("A128").Value = 10
$ws.Range("B128").Value = "Vega Modelo de Temuco"
$ws.Range("C128").Value = "La Araucanía"
$ws.Range("D128").Value = 44504
$ws.Range("E128").Value = 9
$ws.Range("F128").Value = 100112013
$ws.Range("G128").Value = "Alcachofa"
$ws.Range("H128").Value = "Madrigal"
$ws.Range("I128").Value = "Primera"
$ws.Range("J128").Value = 95
$ws.Range("K128").Value = 12000
$ws.Range("L128").Value = 12000
$ws.Range("M128").Value = 12000
$ws.Range("N128").Value = "`$/caja 40 unidades"
$ws.Range("O128").Value = "Región del Maule"
$ws.Range("P128").Value = 300
$ws.Range("Q128").Value = 40
$ws.Range("R128").Value = "Hortaliza"
